$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.17537796497345
$ws.Range("B1").Value = 2.403544664382935
$ws.Range("D1").Value = 2.349050045013428
$ws.Range("E1").Value = 1.207741498947144
